$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 13.5944
$ws.Range("D3").Value = -7.857199999999993
$ws.Range("B4").Value = 5.174600000000003
$ws.Range("C4").Value = -14.3555
$ws.Range("D4").Value = -8.371799999999995
$ws.Range("C5").Value = -14.71650000000001
$ws.Range("E5").Value = 12.42939999999999
$ws.Range("B6").Value = 9.043900000000006
$ws.Range("C6").Value = -11.53100000000001
$ws.Range("B7").Value = 5.240799999999998
$ws.Range("B8").Value = 4.6519
$ws.Range("C8").Value = -11.5438
$ws.Range("D9").Value = -7.577199999999998
$ws.Range("D11").Value = -8.400599999999997
$ws.Range("D14").Value = -6.727699999999993
$ws.Range("B16").Value = 9.116300000000006
$ws.Range("C16").Value = -12.01640000000001
$ws.Range("D18").Value = -8.099999999999993
$ws.Range("B20").Value = 5.102000000000004
$ws.Range("E20").Value = 12.06479999999999
$ws.Range("B21").Value = 4.9817
$ws.Range("C22").Value = -11.08549999999999
$ws.Range("D25").Value = -8.692199999999991
